$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# New "MaxThread" setting row, mirroring the existing Name/Value rows above it
$ws.Range("A5").Value = "MaxThread"
$ws.Range("B5").Value = 4

# Move the selection to B10, matching the saved cursor position
$ws.Activate()
$ws.Range("B10").Select()
